# #272 Ajout d'un scenario de recherche de l'offre d'un professionnel avec un ID Nat PS
# - Bump the "Date" metadata value.
# - Swap the two "Mapping" columns (AK/AL) on the Elements sheet: the
#   "Spécification métier vers l'extension ROR LevelRecourseORSAN" mapping
#   moves before the "RIM Mapping" one (header + data + column widths).

$wb = $excel.ActiveWorkbook

# --- Metadata sheet: refresh the generation Date -------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2024-03-22T16:25:12+00:00"

# --- Elements sheet: swap columns AK (37) and AL (38) --------------------
$wsEl = $wb.Worksheets.Item("Elements")

$colAK = 37
$colAL = 38

# Swap the header row (row 1): "Mapping: RIM Mapping" <-> "Mapping:
# Spécification métier vers l'extension ROR LevelRecourseORSAN"
$headerAK = $wsEl.Cells.Item(1, $colAK).Value2
$headerAL = $wsEl.Cells.Item(1, $colAL).Value2
$wsEl.Cells.Item(1, $colAK).Value = $headerAL
$wsEl.Cells.Item(1, $colAL).Value = $headerAK

# Swap the data rows. Only rows whose AK/AL values actually differ need to
# be touched (rows 2 and 4 hold the same empty value in both columns).
$lastRow = $wsEl.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $akVal = $wsEl.Cells.Item($r, $colAK).Value2
    $alVal = $wsEl.Cells.Item($r, $colAL).Value2
    if ($akVal -cne $alVal) {
        $wsEl.Cells.Item($r, $colAK).Value = $alVal
        $wsEl.Cells.Item($r, $colAL).Value = $akVal
    }
}

# Swap the column widths that go along with the two mapping columns.
$wsEl.Columns.Item($colAK).ColumnWidth = 76.45
$wsEl.Columns.Item($colAL).ColumnWidth = 24.15
